$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C16").Value = "14.04: (PÅSKE)"
$ws.Range("D16").Value = "16.04: (PÅSKE)"

$ws.Range("A17").Value = 17
$ws.Range("B17").Value = "Forberedelse til eksamen"

$ws.Range("B16").Value = "Dataøving 5"

$ws.Range("C17").Value = "21.04: **Speedrun anlyser**"
$ws.Range("D17").Value = "23.04: **Speedrun anlyser**"

[void]$ws.Range("D18").Select()
